# Insert a new data row before the existing row 460 ("Femacal de La Calera" /
# Espinaca weekly price sheet). This pushes the former rows 460..581 down to
# 461..582 and leaves a fresh, blank row 460 to be populated with the new
# weekly observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(460).Insert()

# Populate the newly inserted row 460 with the new observation.
$ws.Cells.Item(460, 1).Value  = 3
$ws.Cells.Item(460, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(460, 3).Value  = "Coquimbo"
$ws.Cells.Item(460, 4).Value  = 45135
$ws.Cells.Item(460, 5).Value  = 5
$ws.Cells.Item(460, 6).Value  = 100112012
$ws.Cells.Item(460, 7).Value  = "Espinaca"
$ws.Cells.Item(460, 8).Value  = "Sin especificar"
$ws.Cells.Item(460, 9).Value  = "Primera"
$ws.Cells.Item(460, 10).Value = 80
$ws.Cells.Item(460, 11).Value = 4000
$ws.Cells.Item(460, 12).Value = 4000
$ws.Cells.Item(460, 13).Value = 4000
$ws.Cells.Item(460, 14).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(460, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(460, 16).Value = 1333
$ws.Cells.Item(460, 17).Value = 3
$ws.Cells.Item(460, 18).Value = "Hortaliza"
